$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the
#     37ca855b... row (row 3) and 624ffa6e... row (row 4) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-13 18:19:03"
$wsOverview.Range("G4").Value = "2016-08-13 18:19:03"

# --- zh-cn sheet: Priority + Handoff/Handback datetimes for rows 3 & 4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-13 18:18:53"
$wsZhCn.Range("H4").Value = "2016-08-13 18:18:53"
$wsZhCn.Range("K3").Value = "2016-08-13 18:19:26"
$wsZhCn.Range("K4").Value = "2016-08-13 18:19:26"

# --- de-de sheet: Handback datetime for rows 3 & 4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-13 18:19:36"
$wsDeDe.Range("K4").Value = "2016-08-13 18:19:36"
